$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The newly appended rows mirror the existing ones: company name, issue
# timestamp, and the comma-formatted share/price figures are stored as plain
# text, while the round number ("회차") is a real number. Pre-format the
# text columns so Excel does not reinterpret "43,315" as a numeric value,
# then strip the temporary formatting back off so the new cells end up
# unstyled, matching the rest of the sheet.
$ws.Range("A10:B23").NumberFormat = "@"
$ws.Range("D10:E23").NumberFormat = "@"

$ws.Cells.Item(10, 1).Value = "에스티팜"
$ws.Cells.Item(10, 2).Value = "2024-09-12 17:20"
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = "43,315"
$ws.Cells.Item(10, 5).Value = "79,648"

$ws.Cells.Item(11, 1).Value = "에스티팜"
$ws.Cells.Item(11, 2).Value = "2024-09-12 17:20"
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 4).Value = "21,971"
$ws.Cells.Item(11, 5).Value = "79,648"

$ws.Cells.Item(12, 1).Value = "에스티팜"
$ws.Cells.Item(12, 2).Value = "2024-09-09 16:50"
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 4).Value = "119,272"
$ws.Cells.Item(12, 5).Value = "79,648"

$ws.Cells.Item(13, 1).Value = "에스티팜"
$ws.Cells.Item(13, 2).Value = "2024-09-04 17:22"
$ws.Cells.Item(13, 3).Value = 2
$ws.Cells.Item(13, 4).Value = "17,577"
$ws.Cells.Item(13, 5).Value = "79,648"

$ws.Cells.Item(14, 1).Value = "에스티팜"
$ws.Cells.Item(14, 2).Value = "2024-09-04 17:22"
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(14, 4).Value = "47,709"
$ws.Cells.Item(14, 5).Value = "79,648"

$ws.Cells.Item(15, 1).Value = "에스티팜"
$ws.Cells.Item(15, 2).Value = "2024-09-03 17:23"
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 4).Value = "32,639"
$ws.Cells.Item(15, 5).Value = "79,648"

$ws.Cells.Item(16, 1).Value = "에스티팜"
$ws.Cells.Item(16, 2).Value = "2024-09-02 17:28"
$ws.Cells.Item(16, 3).Value = 2
$ws.Cells.Item(16, 4).Value = "50,534"
$ws.Cells.Item(16, 5).Value = "79,648"

$ws.Cells.Item(17, 1).Value = "에스티팜"
$ws.Cells.Item(17, 2).Value = "2024-08-29 16:38"
$ws.Cells.Item(17, 3).Value = 2
$ws.Cells.Item(17, 4).Value = "4,394"
$ws.Cells.Item(17, 5).Value = "79,648"

$ws.Cells.Item(18, 1).Value = "에스티팜"
$ws.Cells.Item(18, 2).Value = "2024-08-28 16:47"
$ws.Cells.Item(18, 3).Value = 2
$ws.Cells.Item(18, 4).Value = "43,943"
$ws.Cells.Item(18, 5).Value = "79,648"

$ws.Cells.Item(19, 1).Value = "에스티팜"
$ws.Cells.Item(19, 2).Value = "2024-07-22 17:02"
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = "218,248"
$ws.Cells.Item(19, 5).Value = "68,729"

$ws.Cells.Item(20, 1).Value = "에스티팜"
$ws.Cells.Item(20, 2).Value = "2024-06-14 16:32"
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 4).Value = "145,498"
$ws.Cells.Item(20, 5).Value = "68,729"

$ws.Cells.Item(21, 1).Value = "에스티팜"
$ws.Cells.Item(21, 2).Value = "2024-05-24 17:07"
$ws.Cells.Item(21, 3).Value = 1
$ws.Cells.Item(21, 4).Value = "181,873"
$ws.Cells.Item(21, 5).Value = "68,729"

$ws.Cells.Item(22, 1).Value = "에스티팜"
$ws.Cells.Item(22, 2).Value = "2024-05-17 16:41"
$ws.Cells.Item(22, 3).Value = 1
$ws.Cells.Item(22, 4).Value = "202,243"
$ws.Cells.Item(22, 5).Value = "68,729"

$ws.Cells.Item(23, 1).Value = "에스티팜"
$ws.Cells.Item(23, 2).Value = "2024-04-02 17:22"
$ws.Cells.Item(23, 3).Value = 1
$ws.Cells.Item(23, 4).Value = "125,129"
$ws.Cells.Item(23, 5).Value = "68,729"

$ws.Range("A10:B23").ClearFormats()
$ws.Range("D10:E23").ClearFormats()

